$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.854.49'
$ws.Range('E2').Value = '  -3.35%  '

$ws.Range('D3').Value = '2.912.88'
$ws.Range('E3').Value = '  -4.08%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.59'
$ws.Range('E5').Value = '  -1.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.45'
$ws.Range('E6').Value = '  -6.31%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -2.45%  '

$ws.Range('D9').Value = '2.911.06'
$ws.Range('E9').Value = '  -4.02%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.65'
$ws.Range('E10').Value = '  -2.97%  '

$ws.Range('E11').Value = '  -5.28%  '

$ws.Range('E12').Value = '  -3.96%  '

$ws.Range('E13').Value = '  -3.83%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.39'
$ws.Range('E14').Value = '  -6.80%  '

$ws.Range('E15').Value = '  +1.41%  '

$ws.Range('D16').Value = '3.393.47'
$ws.Range('E16').Value = '  -4.09%  '

$ws.Range('D17').Value = '60.784.85'
$ws.Range('E17').Value = '  -3.36%  '

$ws.Range('E18').Value = '  -5.22%  '

$ws.Range('D19').Value = '2.911.66'
$ws.Range('E19').Value = '  -4.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '428.11'
$ws.Range('E20').Value = '  -5.70%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.54'
$ws.Range('E21').Value = '  -5.14%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.681'
$ws.Range('E22').Value = '  -2.47%  '

$ws.Range('E23').Value = '  -5.70%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.62'
$ws.Range('E24').Value = '  -3.02%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.23'
$ws.Range('E25').Value = '  -3.07%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.75'
$ws.Range('E26').Value = '  -3.79%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.93'
$ws.Range('E27').Value = '  -4.23%  '

$ws.Range('E28').Value = '  +0.00%  '

$ws.Range('E29').Value = '  +0.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.21'
$ws.Range('E30').Value = '  -3.63%  '

$ws.Range('E31').Value = '  -3.43%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.17'
$ws.Range('E32').Value = '  -4.06%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.49'
$ws.Range('E33').Value = '  -4.08%  '

$ws.Range('D35').Value = '0.0₃0874'
$ws.Range('E35').Value = '  +1.58%  '

$ws.Range('E36').Value = '  -3.11%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.61'
$ws.Range('E37').Value = '  -5.52%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.01'
$ws.Range('E38').Value = '  -6.29%  '

$ws.Range('E41').Value = '  -5.49%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.59'
$ws.Range('E42').Value = '  -6.08%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.297'
$ws.Range('E43').Value = '  -3.24%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.28'
$ws.Range('E44').Value = '  -6.11%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '378.55'
$ws.Range('E45').Value = '  -3.43%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0350'
$ws.Range('E46').Value = '  -3.22%  '

$ws.Range('D47').Value = '2.688.74'
$ws.Range('E47').Value = '  -1.27%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.26'
$ws.Range('E48').Value = '  -0.90%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.41'
$ws.Range('E50').Value = '  -1.85%  '

$ws.Range('E51').Value = '  -2.62%  '

# Row 39 and 40: Kaspa/OKB swap positions with updated data
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.126'
$ws.Range('E39').Value = '  -4.18%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.45'
$ws.Range('E40').Value = '  -1.93%  '
